# Insert a new column before the existing data so the old column A
# (segment names) becomes column B, old B (PercActivations) becomes C,
# and old C (PercSegmentAreas) becomes D. Then populate the new column A
# with a numeric index (0..18) for each segment row, and add the
# "segments" header above the (now) segment-name column B.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift everything one column to the right by inserting a new column A.
$ws.Range("A1").EntireColumn.Insert()

# Give the new header cell (B1) the same formatting as the other header
# cells (bold / bordered / centered style used in row 1), then set its
# text.
$ws.Range("C1").Copy()
$ws.Range("B1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("B1").Value = "segments"

# Give the new index cells (A2:A20) the same formatting as the other
# label-column cells (B2:B20, which carry the original header-row style),
# then fill them with the 0-based segment index.
$ws.Range("B2").Copy()
$ws.Range("A2:A20").PasteSpecial(-4122)  # xlPasteFormats

for ($i = 0; $i -le 18; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $i
}

$excel.CutCopyMode = $false
